$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new quarter-end dates in V1:X1, copying the formatting of U1 ---
$ws.Range("U1").Copy()
$ws.Range("V1:X1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(1, 22).Value = "31/12/2023"
$ws.Cells.Item(1, 23).Value = "31/03/2024"
$ws.Cells.Item(1, 24).Value = "30/06/2024"

# --- Data rows: values concatenated from the 31/12/2023, 31/03/2024 and 30/06/2024 balance sheets ---
$newQuarters = @{
    2 = @(4604387.84, 4660456.96, 4087429.888)
    3 = @(1950392.96, 1997640.96, 1421421.952)
    4 = @(1188222.976, 1126301.056, 565830.976)
    5 = @(0, 0, 0)
    6 = @(707555.968, 802556.032, 771752)
    7 = @(0, 0, 0)
    8 = @(0, 0, 0)
    9 = @(27299, 27553, 41818)
    10 = @(0, 0, 0)
    11 = @(27315, 41231, 42021)
    12 = @(208206, 227638, 232198)
    13 = @(0, 0, 0)
    14 = @(0, 0, 0)
    15 = @(0, 0, 0)
    16 = @(0, 0, 0)
    17 = @(0, 0, 0)
    18 = @(0, 0, 0)
    19 = @(182807.008, 200406, 203156.992)
    20 = @(0, 0, 0)
    21 = @(0, 0, 0)
    22 = @(0, 0, 0)
    23 = @(157640.992, 150374, 145647.008)
    24 = @(2288147.968, 2284804.096, 2288163.072)
    25 = @(0, 0, 0)
    26 = @(4604387.84, 4660456.96, 4087429.888)
    27 = @(1441553.024, 1479239.04, 1067249.024)
    28 = @(108366, 120592, 100760)
    29 = @(45575, 34110, 28624)
    30 = @(17611, 17028, 22367)
    31 = @(13754, 12830, 12842)
    32 = @(0, 0, 0)
    33 = @(5, 5, 5)
    34 = @(1256242.04, 1294674.04, 902651)
    35 = @(0, 0, 0)
    36 = @(0, 0, 0)
    37 = @(280688.992, 281388.992, 128477)
    38 = @(63675, 61866, 60319)
    39 = @(0, 0, 0)
    40 = @(204844, 207940, 55272)
    41 = @(6002, 5733, 6805)
    42 = @(0, 0, 0)
    43 = @(6168, 5850, 6081)
    44 = @(0, 0, 0)
    45 = @(0, 0, 0)
    46 = @(0, 0, 0)
    47 = @(2882146.048, 2899828.992, 2891704.064)
    48 = @(2940141.056, 2940141.056, 2868290.048)
    49 = @(13856, 7175, -21136)
    50 = @(0, 0, 0)
    51 = @(0, 0, 0)
    52 = @(-71851, -47381, 42792)
    53 = @(0, -106, 1758)
    54 = @(0, 0, 0)
    55 = @(0, 0, 0)
    56 = @(0, 0, 0)
    59 = @(345918.976, 320612, 335953.984)
    60 = @(-181529.008, -168076.992, -175919.008)
    61 = @(164390.016, 152535.008, 160035.008)
    62 = @(-66182.992, -67072, -73507)
    63 = @(-68000, -59896, -63177)
    64 = @(-3887, -2626, -2090)
    65 = @(4675, 667, 4499)
    66 = @(0, 0, 0)
    67 = @(0, 0, 0)
    68 = @(-111601, 9251, -4214)
    69 = @(37010.008, 34687, 26044)
    70 = @(-148611.008, -25436, -30258)
    74 = @(-80606, 32859, 21546)
    75 = @(-7353, -26203, -5862)
    76 = @(41952, 17814, 2638)
    79 = @(0, 0, 0)
    80 = @(-46006, 24470, 18322)
}
foreach ($row in $newQuarters.Keys) {
    $vals = $newQuarters[$row]
    $ws.Cells.Item($row, 22).Value = $vals[0]
    $ws.Cells.Item($row, 23).Value = $vals[1]
    $ws.Cells.Item($row, 24).Value = $vals[2]
}

# --- Blank separator rows: mirror the (empty) V/W/X cells alongside the existing blank columns ---
$blankRows = @(57, 58, 71, 72, 73, 77, 78)
foreach ($row in $blankRows) {
    $ws.Cells.Item($row, 22).Value = ""
    $ws.Cells.Item($row, 23).Value = ""
    $ws.Cells.Item($row, 24).Value = ""
}

$ws.Range("A1:X80").Select()
